$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "Arjun"
$ws.Range("E7").Value = "saivamsi"

$ws.Range("E2").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("E7").Select()
